# Changes of 6th May 2022
# CheetahProcessing.xlsx test data refresh: rows 2-22 get new
# ShipmentTrackNum (column C) / PackageTrackNum (column D) values.
#
# The source cells are plain shared-string text (12-digit numbers stored
# as Text, not Number), so every new value is written with a leading
# apostrophe to force Excel to keep it as text instead of re-inferring a
# numeric type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  Value = "320018483205"; HasD = $false },
    @{ Row = 3;  Value = "320018483238"; HasD = $false },
    @{ Row = 4;  Value = "320018483260"; HasD = $false },
    @{ Row = 5;  Value = "320018483282"; HasD = $true  },
    @{ Row = 6;  Value = "320018483330"; HasD = $true  },
    @{ Row = 7;  Value = "320018483352"; HasD = $true  },
    @{ Row = 8;  Value = "320018483385"; HasD = $false },
    @{ Row = 9;  Value = "320018483411"; HasD = $false },
    @{ Row = 10; Value = "320018483444"; HasD = $false },
    @{ Row = 11; Value = "320018483466"; HasD = $false },
    @{ Row = 12; Value = "320018483503"; HasD = $false },
    @{ Row = 13; Value = "320018475505"; HasD = $true  },
    @{ Row = 14; Value = "320018475538"; HasD = $true  },
    @{ Row = 15; Value = "320018475550"; HasD = $true  },
    @{ Row = 16; Value = "320018475582"; HasD = $true  },
    @{ Row = 17; Value = "320018475696"; HasD = $true  },
    @{ Row = 18; Value = "320018475733"; HasD = $false },
    @{ Row = 19; Value = "320018475766"; HasD = $false },
    @{ Row = 20; Value = "320018475799"; HasD = $false },
    @{ Row = 21; Value = "320018475836"; HasD = $false },
    @{ Row = 22; Value = "320018475869"; HasD = $false }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $text = $r.Value

    $ws.Cells.Item($rowNum, 3).Value = "'" + $text
    if ($r.HasD) {
        $ws.Cells.Item($rowNum, 4).Value = "'" + $text
    }
}
